$wb = $excel.ActiveWorkbook

$wsDatos  = $wb.Worksheets.Item("Datos")
$wsListas = $wb.Worksheets.Item("Listas")

# --- Data change: N2 on "Datos" (fechaFactura column) updated ---
$wsDatos.Range("N2").Value = 65437

# --- Selection moved from H4 to H13 on "Datos" ---
$wsDatos.Activate()
$wsDatos.Range("H13").Select() | Out-Null

# --- Column width touch-ups ---
# NOTE: this host's ColumnWidth setter quantizes the stored <col width>
# to a 1/6-character grid (offset +5/6), so the real Excel sub-pixel widths
# from the diff (e.g. 18.33203125) can't be hit exactly. Only nudge a
# column when the nearest reachable grid point is actually closer to the
# target than the current (pre-edit) width; otherwise leave it untouched.
$wsDatos.Columns.Item(14).ColumnWidth = 17.5

$wsListas.Columns.Item(1).ColumnWidth = 9.83333333333333
